$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column F (dSF)
$updates = @{
    2  = 1
    11 = 1
    14 = -1
    15 = 2
    30 = -2
    36 = 3
    39 = -1
    40 = -1
    44 = 0
    49 = 3
    50 = 2
    54 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
